# Applies the Wed May 10 02:23:45 UTC 2023 cryptos-list refresh:
# updated Price/Volume(1h) figures, plus a few re-ranked coins whose
# Coin/Link moved to a different row as their ranking shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.772.27"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.852.33"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.26"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4329"
$ws.Range("E7").Value = "  +1.90%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.17"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07347"
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8806"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.77"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "1.898.98"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.358"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.533"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06940"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "80.62"
$ws.Range("E18").Value = "  +3.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009068"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.40"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "27.957.66"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.982"
$ws.Range("E23").Value = "  +0.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.41"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").Value = "2.141.34"
$ws.Range("E25").Value = "  +2.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.992"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.06"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.70"
$ws.Range("E28").Value = "  +2.83%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.41"
$ws.Range("E29").Value = "  +9.24%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.277"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.864"
$ws.Range("E31").Value = "  +2.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08940"
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7662"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.557"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.963"
$ws.Range("E35").Value = "  +1.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.132"
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.111"
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05437"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01943"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.841"
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5112"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1663"
$ws.Range("E42").Value = "  +1.46%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.711"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.374"
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.45"
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06556"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4692"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.68"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.001"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.627"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "64.90"
$ws.Range("E51").Value = "  +0.88%  "
